$wb = $excel.ActiveWorkbook
$itemWs = $wb.Worksheets.Item("Item")
$roundWs = $wb.Worksheets.Item("Round")

# Fix the misspelled item name "GRANADE LAUNCHER" -> "GRENADE LAUNCHER"
$itemWs.Range("B10").Value = "GRENADE LAUNCHER"

# Fill in the item/tower descriptions
$itemWs.Range("D3").Value = 'A drive-by favorite, the uzi is an affordable, fuuly automatic solution to all your semi-automatic problems.'
$itemWs.Range("D4").Value = 'Rednecks and zombie slayers agree, shotguns are always a good decision. Don’t forget… shop smart, shop S-mart.'
$itemWs.Range("D5").Value = 'You’re standard automatic rifle. Accurate and powerful. A step above that cute lil’ uzi.'
$itemWs.Range("D6").Value = 'A ballsier, automatic version of the combat shotgun. For when you really need to call down the rain.'
$itemWs.Range("D7").Value = 'Perfect for getting rid of termites, hilarious WWII reenactments with friends, or even when your son’s little league game runs long. The flame thrower. Don’t leave home without it.'
$itemWs.Range("D8").Value = 'When Sylvester stalone needs to cut down row after row of Vietcong soldiers and emerge with nothing but stainless steel abs and the sweat on his brow, you better believe he brings his machine gun.'
$itemWs.Range("D9").Value = '“There’s a chainsaw? Sweet!” - You'
$itemWs.Range("D10").Value = 'Perfect for launching grenades.'
$itemWs.Range("D11").Value = 'This baby focuses a satellite mounted laser at your feeble enemies and promptly deatomizes them.'
$itemWs.Range("D12").Value = 'Handheld tele-geo-dynamics manipulating oscillator generator capable of knocking down troops and damaging vehicles.'
$itemWs.Range("D13").Value = 'A gun so big it overlaps the buy button.'
$itemWs.Range("D15").Value = 'This giant fan has a slowing effect that will stop those puny swordmen right in their trakcs. LOL!'
$itemWs.Range("D16").Value = 'This beastbox initiates minor tremors causing foot soldiers to fall over and armored units to take damage.'
$itemWs.Range("D17").Value = 'Known by scientists as the “Magnificant Wallopping Van De Graaff Machine”, the tesla tower delivers lethally concentrated bolts of electricity at a short range.'
$itemWs.Range("D18").Value = 'Fires heat seeking missiles that damage all units within it’s blast radius.'
$itemWs.Range("D19").Value = 'The M6 Laser Tower fires a very powerful blast of energy at the toughest enemy on the screen.'
$itemWs.Range("D20").Value = 'Flamer Joe likes to set people on fire.'
$itemWs.Range("D21").Value = 'Generates a protective shield that cuts all damage taken by a percentage.'

# Make the Item sheet the active tab/selection, as it was when the
# descriptions + icon info were being edited
$itemWs.Activate()
$itemWs.Range("G14").Select()

# Keep the Round sheet selection where it was left
$roundWs.Range("H25").Select()
